$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Logs" -------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

# Add the new log row (row 13) with the new testmail entry
$logs.Range("A13").Value = "Kun jij dit even regelen?"
$logs.Range("B13").Value = "mailmind.test@zohomail.eu"
$logs.Range("C13").Value = "Testmail #1: Kun jij dit even regelen?"
$logs.Range("D13").Value = "Planning / Afspraak"
$logs.Range("E13").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$logs.Range("F13").Value = "2025-08-04 20:26:35"
$logs.Range("G13").Value = "Ja"
$logs.Range("H13").Value = "Ja"
$logs.Range("I13").Value = "Nee"
$logs.Range("J13").Value = "Nee"

# Extend the existing conditional-formatting rules (D, G, H, I, J columns)
# down to the new row 13, without altering the rules themselves.
function Extend-CFRange($ws, $oldRange, $newRange) {
    $fcs = $ws.Range($oldRange).FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($ws.Range($newRange))
    }
}

Extend-CFRange $logs "D2:D12" "D2:D13"
Extend-CFRange $logs "G2:G12" "G2:G13"
Extend-CFRange $logs "H2:H12" "H2:H13"
Extend-CFRange $logs "I2:I12" "I2:I13"
Extend-CFRange $logs "J2:J12" "J2:J13"

# --- Sheet 2: "Dashboard" --------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

# Swap the category summary rows: row2 becomes "Planning / Afspraak" (4),
# row3 becomes "Opvolging / Status" (3)
$dash.Range("A2").Value = "Planning / Afspraak"
$dash.Range("B2").Value = 4
$dash.Range("A3").Value = "Opvolging / Status"
$dash.Range("B3").Value = 3
